# Weekly fruit/vegetable price update:
# Insert a new week's data row at row 3 (before the old row 3), shifting
# all existing data rows (old 3..9) down to (4..10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 3; this pushes rows 3-9 down to 4-10
# and extends the used range / dimension accordingly.
$ws.Rows("3").Insert()

# Fill in the new row 3 with this week's data.
$ws.Cells.Item(3, 1).Value = 11
$ws.Cells.Item(3, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(3, 3).Value = "Bíobío"
$ws.Cells.Item(3, 4).Value = 44526
$ws.Cells.Item(3, 5).Value = 8
$ws.Cells.Item(3, 6).Value = 300000000
$ws.Cells.Item(3, 7).Value = "Espárragos"
$ws.Cells.Item(3, 8).Value = "Sin especificar"
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 100
$ws.Cells.Item(3, 11).Value = 1500
$ws.Cells.Item(3, 12).Value = 1600
$ws.Cells.Item(3, 13).Value = 1550
$ws.Cells.Item(3, 14).Value = "$/kilo"
$ws.Cells.Item(3, 15).Value = "Provincia de Linares"
$ws.Cells.Item(3, 16).Value = 1550
$ws.Cells.Item(3, 17).Value = 1
$ws.Cells.Item(3, 18).Value = "Hortaliza"
